# fix: Correct Problem Reference
#
# In the solution to "سوال ۸" ("Question 8") the text refers back to
# "سوال ۷" ("question 7") but is missing the comma that should follow
# the reference, i.e. "... در سوال ۷ در صورتی ..." should read
# "... در سوال ۷، در صورتی ...".
#
# The Persian/Arabic comma "،" needs to be inserted as its own run
# (matching the formatting of the "۷" run) right after the existing
# "۷" run and before the "_GoBack" bookmark that immediately follows
# it, without disturbing any of the surrounding runs/bookmark.

$d = $word.ActiveDocument

# Locate the unique phrase "سوال ۸: در سوال ۷" - this pins down exactly
# the "۷" that is the mis-punctuated back-reference to question 7
# (there are other "۷" occurrences in the document, but only this one
# is immediately preceded by "سوال ۸: در سوال ").
$rFind = $d.Content
$found = $rFind.Find.Execute("سوال ۸: در سوال ۷", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target 'سوال ۸: در سوال ۷' reference."
}

# Character position immediately after the "۷" (i.e. where the bookmark
# currently starts).
$posAfterDigit = $rFind.End

# Insert the comma right after the found text. (This also happens to
# normalize/merge the neighbouring same-formatted runs in this
# paragraph - we fix that up below so the "۷" and "،" end up as their
# own distinct runs, matching the original document's run layout.)
$rFind.Collapse(0)
$rFind.InsertAfter("،")

# Re-assert the run boundaries: touching each character's formatting
# (a harmless Bold on/off round trip) forces the engine to keep it as
# its own run instead of silently merging it into its neighbour.
$rDigit = $d.Range($posAfterDigit - 1, $posAfterDigit)
$rDigit.Font.Bold = 1
$rDigit.Font.Bold = 0

$rComma = $d.Range($posAfterDigit, $posAfterDigit + 1)
$rComma.Font.Bold = 1
$rComma.Font.Bold = 0

Write-Output "Inserted Persian comma after the 'سوال ۷' reference."
